$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.516.80"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "2.487.37"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "314.01"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "94.58"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "33.64"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("D13").Value = "7.00"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "2.874.21"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "2.478.35"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "41.491.59"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "11.30"
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").Value = "68.96"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "237.50"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "24.22"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "9.79"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "36.71"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").Value = "152.46"
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("D32").Value = "5.51"
$ws.Range("E32").Value = "  -5.92%  "
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("D34").Value = "18.15"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("D35").Value = "0.0759"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "3.09"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -10.47%  "
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -5.11%  "
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "19.87"
$ws.Range("E43").Value = "  -8.55%  "
$ws.Range("D44").Value = "1.995.27"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "0.0287"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "3.03"
$ws.Range("E46").Value = "  -6.81%  "
$ws.Range("D47").Value = "8.87"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("D48").Value = "2.735.19"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "70.11"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "97.15"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.179"
$ws.Range("E51").Value = "  -5.19%  "
